$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

$ws.Cells.Item(10, 1).Value = 43756
$ws.Cells.Item(10, 1).NumberFormat = "d-mmm-yy"
$ws.Cells.Item(10, 2).Value = 43756
$ws.Cells.Item(10, 2).NumberFormat = "d-mmm-yy"
$ws.Cells.Item(1, 5).Value = "tuple_ints"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(3, 5).Value = $false
$ws.Cells.Item(4, 5).Value = $true
$ws.Cells.Item(6, 5).Value = "hello"
$ws.Cells.Item(8, 5).Value = "1, 2, 3"
$ws.Cells.Item(9, 5).Value = "20 8"
$ws.Cells.Item(10, 5).Value = 43756
$ws.Cells.Item(10, 5).NumberFormat = "d-mmm-yy"
$ws.Cells.Item(5, 5).Value = "19twenty3"
$ws.Cells.Item(1, 6).Value = "boolean"
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(3, 6).Value = $false
$ws.Cells.Item(4, 6).Value = $true
$ws.Cells.Item(6, 6).Value = "hello"
$ws.Cells.Item(8, 6).Value = "1, 2, 3"
$ws.Cells.Item(9, 6).Value = "20 8"
$ws.Cells.Item(10, 6).Value = 43756
$ws.Cells.Item(10, 6).NumberFormat = "d-mmm-yy"
$ws.Cells.Item(5, 6).Value = "19twenty3"
$ws.Cells.Item(1, 8).Value = "ints_with_missing"
$ws.Cells.Item(2, 8).Value = 1
$ws.Cells.Item(3, 8).Value = $false
$ws.Cells.Item(4, 8).Value = $true
$ws.Cells.Item(6, 8).Value = "hello"
$ws.Cells.Item(8, 8).Value = 25.5
$ws.Cells.Item(9, 8).Value = 858
$ws.Cells.Item(10, 8).Value = 43756
$ws.Cells.Item(10, 8).NumberFormat = "d-mmm-yy"
$ws.Cells.Item(5, 8).Value = "19twenty3"
$ws.Cells.Item(1, 7).Value = "all_numbers"
$ws.Cells.Item(2, 7).Value = 1
$ws.Cells.Item(3, 7).Value = 2
$ws.Cells.Item(4, 7).Value = 3
$ws.Cells.Item(6, 7).Value = 77557357
$ws.Cells.Item(8, 7).Value = 25.5
$ws.Cells.Item(9, 7).Value = 858
$ws.Cells.Item(10, 7).Value = 23
$ws.Cells.Item(5, 7).Value = 465
$ws.Cells.Item(1, 9).Value = "float"
$ws.Cells.Item(2, 9).Value = 1
$ws.Cells.Item(3, 9).Value = $false
$ws.Cells.Item(4, 9).Value = $true
$ws.Cells.Item(6, 9).Value = "hello"
$ws.Cells.Item(8, 9).Value = 25.5
$ws.Cells.Item(9, 9).Value = 858
$ws.Cells.Item(10, 9).Value = 43756
$ws.Cells.Item(10, 9).NumberFormat = "d-mmm-yy"
$ws.Cells.Item(5, 9).Value = "19twenty3"
$ws.Cells.Item(1, 4).Value = "first_digit_missing"
$ws.Cells.Item(2, 4).Value = 2
$ws.Cells.Item(3, 4).Value = $false
$ws.Cells.Item(4, 4).Value = $true
$ws.Cells.Item(6, 4).Value = "hello"
$ws.Cells.Item(8, 4).Value = "helper 30"
$ws.Cells.Item(9, 4).Value = "20 manager"
$ws.Cells.Item(10, 4).Value = 43756
$ws.Cells.Item(10, 4).NumberFormat = "d-mmm-yy"
$ws.Cells.Item(5, 4).Value = "'1"
$ws.Cells.Item(7, 3).Value = 57
$ws.Cells.Item(7, 7).Value = 12541
$ws.Cells.Item(3, 3).Value = 457
$ws.Cells.Item(4, 3).Value = 63.5
$ws.Cells.Item(6, 3).Value = 3
$ws.Cells.Item(10, 3).Value = "30 manager"
$ws.Cells.Item(5, 3).Value = "'4"

$ws.Columns.Item(3).ColumnWidth = 11.140625
$ws.Columns.Item(4).ColumnWidth = 17.5703125
$ws.Columns.Item(5).ColumnWidth = 10.28515625
$ws.Columns.Item(6).ColumnWidth = 10.28515625
$ws.Columns.Item(7).ColumnWidth = 12
$ws.Columns.Item(8).ColumnWidth = 17.28515625
$ws.Columns.Item(9).ColumnWidth = 10.28515625
$ws.Range("C5").Select()
